$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with fresh data pulled on
# Sun Feb  5 17:46:38 UTC 2023. Values are plain text (not numeric/percent)
# in the source sheet, so force text format before assigning to avoid Excel
# auto-converting them to numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.23%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.50%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.428"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.13%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08118"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.63%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.717"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.93%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.328"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.41%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.893"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.37%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9447"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.14%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1182"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1889"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.90%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09672"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.96%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04230"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001293"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.28%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006009"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.58%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.558"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.41%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.53%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.799"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.32%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1361"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.08%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2608"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.05%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04393"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.46%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001243"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.36%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004323"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.44%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001241"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.73%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004020"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "32.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02662"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.67%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.07%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007842"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.39%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009779"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.07%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1401"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002129"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.46%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009617"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-13.27%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007342"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.77%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000756"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.76%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003466"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.02%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002288"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.42%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002117"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.76%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002016"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.76%"
